# RAWDOCKET project 4 - add STLC (St. Louis City) and STLCC docket rows.
# Appends three new case rows (53-55) to the single worksheet, matching the
# same column layout used by every existing row (A:W).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 53 - Buchanan / Washington University Physician Network
# ---------------------------------------------------------------------------
$ws.Range("A53").Value = "Buchanan"
$ws.Range("B53").Value = "Shelly"
$ws.Range("C53").Value = "Shelly K. Buchanan"
$ws.Range("D53").Value = ""
$ws.Range("E53").Value = "Washington University Physician Network"
$ws.Range("F53").Value = ""
$ws.Range("G53").Value = ""
$ws.Range("H53").Value = "St. Louis City"
$ws.Range("I53").Value = "MO"
$ws.Range("J53").Value = ""
$ws.Range("K53").Value = "2022-AC03429"
$ws.Range("L53").Value = "3375223"
$ws.Range("M53").Value = "1/13/2021"
$ws.Range("N53").Value = "10:30 AM"
$ws.Range("O53").Value = "Docket"
$ws.Range("P53").Value = "1/2/2021"
$ws.Range("Q53").Value = "No Service"
$ws.Range("R53").Value = "10/3/2020"
$ws.Range("S53").Value = "131"
$ws.Range("T53").Value = ""
$ws.Range("U53").Value = ""
$ws.Range("V53").Value = "CCM"
$ws.Range("W53").Value = 859.64

# ---------------------------------------------------------------------------
# Row 54 - Coleman / Washington University Physician Network
# ---------------------------------------------------------------------------
$ws.Range("A54").Value = "Coleman"
$ws.Range("B54").Value = "Danielle"
$ws.Range("C54").Value = "Danielle P. Coleman"
$ws.Range("D54").Value = ""
$ws.Range("E54").Value = "Washington University Physician Network"
$ws.Range("F54").Value = ""
$ws.Range("G54").Value = ""
$ws.Range("H54").Value = "St. Louis City"
$ws.Range("I54").Value = "MO"
$ws.Range("J54").Value = ""
$ws.Range("K54").Value = "2022-AC03533"
$ws.Range("L54").Value = "3375504"
$ws.Range("M54").Value = "1/13/2021"
$ws.Range("N54").Value = "10:30 AM"
$ws.Range("O54").Value = "Docket"
$ws.Range("P54").Value = "10/4/2020"
$ws.Range("Q54").Value = "Personal"
$ws.Range("S54").Value = "101"
$ws.Range("T54").Value = ""
$ws.Range("U54").Value = ""
$ws.Range("V54").Value = "CCM"
$ws.Range("W54").Value = 754.18
$ws.Rows.Item(54).RowHeight = 14.25

# ---------------------------------------------------------------------------
# Row 55 - Ewing / Southfield Partnership, L.P. d/b/a Southfield Apartments
# ---------------------------------------------------------------------------
$ws.Range("A55").Value = "Ewing"
$ws.Range("B55").Value = "Jamina"
$ws.Range("C55").Value = "Jamina J. Ewing"
$ws.Range("D55").Value = ""
$ws.Range("E55").Value = "Southfield Partnership, L.P. d/b/a Southfield Apartments"
$ws.Range("F55").Value = ""
$ws.Range("G55").Value = ""
$ws.Range("H55").Value = "St. Louis City"
$ws.Range("I55").Value = "MO"
$ws.Range("J55").Value = ""
$ws.Range("K55").Value = "2022-AC01514"
$ws.Range("L55").Value = "3375941"
$ws.Range("M55").Value = "1/13/2021"
$ws.Range("N55").Value = "10:30 AM"
$ws.Range("O55").Value = "Docket"
$ws.Range("P55").Value = "12/10/2020"
$ws.Range("Q55").Value = "Personal"
$ws.Range("S55").Value = "101"
$ws.Range("T55").Value = ""
$ws.Range("U55").Value = ""
$ws.Range("V55").Value = "CCM"
$ws.Range("W55").Value = 1813.45

# ---------------------------------------------------------------------------
# Update the worksheet's view to match the freshly-entered rows: scroll down
# towards the new block and leave the 3 newly-typed rows selected, same as
# Excel leaves behind after typing/pasting a new block of rows.
# ---------------------------------------------------------------------------
$excel.Goto($ws.Range("A46"), $true)
$ws.Range("A53:XFD55").Select()
